$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# New header labels in columns A-D (row 1)
$ws.Range("C1").Value = "2theta"
$ws.Range("A1").Value = "d2theta"
$ws.Range("D1").Value = "theta_rad"
$ws.Range("B1").Value = "d2theta_rad"

# Update the active selection to N4, as in the final workbook state
$ws.Range("N4").Select()

# Reposition/resize the first two chart objects (the user nudged them while
# finishing the report layout)
$chart1 = $ws.ChartObjects(1)
$chart1.Left = 205.232421875
$chart1.Top = 115
$chart1.Width = 446.5
$chart1.Height = 216.5

$chart2 = $ws.ChartObjects(2)
$chart2.Left = 669.232421875
$chart2.Top = 115.5
$chart2.Width = 412.98828125
$chart2.Height = 216.5
